$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("orders")

# Order 1 (row 2): fix status column (was mistakenly in D, belongs in C),
# refresh the order timestamp, and record the finished order contents.
$ws.Range("D2").ClearContents()
$ws.Range("C2").Value = "Создается"
$ws.Range("F2").Value = "2022-08-21 15:37:32.832707"
$ws.Range("G2").Value = "мороженое, , чай, лимонад, смузи, смузи, мороженое, вафлю, молочный_коктель"

# Order 2 (row 3): same status-column fix, new timestamp and order contents.
$ws.Range("D3").ClearContents()
$ws.Range("C3").Value = "Создается"
$ws.Range("F3").Value = "2022-08-21 15:38:51.982234"
$ws.Range("G3").Value = "-чай -чай -чай -смузи -мороженое -лимонад -вафлю -молочный_коктель -чай"

# Orders 3-5 (rows 4-6) are removed / not created anymore.
$ws.Range("A4:G6").ClearContents()

# Drop the leftover per-row style override left behind by the earlier rows,
# matching the rest of the (unstyled) data rows.
$ws.Range("A2:G22").EntireRow.ClearFormats()
